$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mojibake characters in the Regional Economic Communities footnote (shared string)
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development""."

# Update recalculated statistical values in rows 62-98
$ws.Range("G62").Value = 48.279545454545499
$ws.Range("C63").Value = 2.3294736842105301
$ws.Range("D63").Value = 8.6094736842105295
$ws.Range("E63").Value = 23.057894736842101
$ws.Range("F63").Value = 34.991578947368403
$ws.Range("G63").Value = 42.591578947368397
$ws.Range("C65").Value = 4.95
$ws.Range("D65").Value = 21.774999999999999
$ws.Range("E65").Value = 52.03
$ws.Range("F65").Value = 35.174999999999997
$ws.Range("G65").Value = 43.274999999999999
$ws.Range("C66").Value = 10.0407142857143
$ws.Range("D66").Value = 21.925714285714299
$ws.Range("E66").Value = 39.810714285714297
$ws.Range("F66").Value = 37.119285714285702
$ws.Range("G66").Value = 44.3920863309352
$ws.Range("G68").Value = 45.038095238095202
$ws.Range("C76").Value = 2.1857142857142899
$ws.Range("D76").Value = 13.5285714285714
$ws.Range("E76").Value = 41.242857142857098
$ws.Range("F76").Value = 37.314285714285703
$ws.Range("G76").Value = 44.8857142857143
$ws.Range("D82").Value = 51.130769230769197
$ws.Range("C83").Value = 1.9747126436781599
$ws.Range("D83").Value = 7.6977011494252903
$ws.Range("E83").Value = 21.778160919540198
$ws.Range("F83").Value = 35.052873563218398
$ws.Range("G83").Value = 42.605747126436803
$ws.Range("E86").Value = 72.534999999999997
$ws.Range("C87").Value = 5.4249999999999998
$ws.Range("D87").Value = 22.4166666666667
$ws.Range("E87").Value = 54.85
$ws.Range("F87").Value = 35.908333333333303
$ws.Range("G87").Value = 43.683333333333302
$ws.Range("F89").Value = 38.133333333333297
$ws.Range("E91").Value = 85.757142857142895
$ws.Range("D95").Value = 68.193333333333399
$ws.Range("C97").Value = 29.669230769230801
$ws.Range("D97").Value = 57.492307692307698
$ws.Range("E97").Value = 83.796153846153899
$ws.Range("F97").Value = 41.719230769230798
$ws.Range("G97").Value = 48.257692307692302
$ws.Range("H97").Value = 30.152307692307701
$ws.Range("I97").Value = 31.583846153846199
$ws.Range("J97").Value = 28.876923076923099
$ws.Range("K97").Value = 35.085769230769202
$ws.Range("L97").Value = 33.602307692307697
$ws.Range("M97").Value = 36.261153846153903
$ws.Range("N97").Value = 9980.3927440706302
$ws.Range("C98").Value = 8.51
$ws.Range("D98").Value = 31.11
$ws.Range("E98").Value = 64.510000000000005
$ws.Range("F98").Value = 37.78
$ws.Range("G98").Value = 45.28
$ws.Range("H98").Value = 9.44
$ws.Range("I98").Value = 9.2706666666666706
$ws.Range("J98").Value = 9.5
$ws.Range("K98").Value = 11.736000000000001
$ws.Range("L98").Value = 10.93
$ws.Range("M98").Value = 12.135999999999999
$ws.Range("N98").Value = 20083.719205813999
